$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1928
$ws.Range("F3").Value = 1538
$ws.Range("F4").Value = 908
$ws.Range("F5").Value = 801
$ws.Range("F6").Value = 13469
$ws.Range("F7").Value = 13310
$ws.Range("F8").Value = 1027
$ws.Range("F9").Value = 786
$ws.Range("F11").Value = 574
$ws.Range("F12").Value = 73
$ws.Range("F15").Value = 700
$ws.Range("F16").Value = 2108
$ws.Range("F17").Value = 27
$ws.Range("F18").Value = 70
$ws.Range("F19").Value = 50
$ws.Range("F25").Value = 446
$ws.Range("F26").Value = 778
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 138
$ws.Range("F8").Value = 613
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 66
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1928
$ws.Range("F4").Value = 1538
$ws.Range("F5").Value = 908
$ws.Range("F7").Value = 801
$ws.Range("F8").Value = 13469
$ws.Range("F9").Value = 13310
$ws.Range("F10").Value = 1027
$ws.Range("F11").Value = 786
$ws.Range("F13").Value = 574
$ws.Range("F14").Value = 73
$ws.Range("F17").Value = 700
$ws.Range("F20").Value = 2108
$ws.Range("F21").Value = 27
$ws.Range("F22").Value = 70
$ws.Range("F23").Value = 50
$ws.Range("F28").Value = 66
$ws.Range("F32").Value = 446
$ws.Range("F33").Value = 778
$ws.Range("F34").Value = 138
$ws.Range("F35").Value = 613
